{"js": "// Replace the date line and the 25 two-digit-by-two-digit multiplication\n// answers in the practice sheet with the new day's values. Every old\n// value below is unique in the document, so a simple exact-text search\n// and replace (one pair at a time) reproduces the diff exactly.\nconst replacements = [\n  [\"2025-03-27 Thursday\", \"2025-03-28 Friday\"],\n  [\"56\u00d749=2744\", \"71\u00d754=3834\"],\n  [\"53\u00d737=1961\", \"77\u00d746=3542\"],\n  [\"13\u00d742=546\", \"76\u00d788=6688\"],\n  [\"71\u00d770=4970\", \"96\u00d738=3648\"],\n  [\"96\u00d719=1824\", \"30\u00d764=1920\"],\n  [\"60\u00d777=4620\", \"58\u00d731=1798\"],\n  [\"19\u00d782=1558\", \"66\u00d732=2112\"],\n  [\"46\u00d794=4324\", \"12\u00d788=1056\"],\n  [\"42\u00d788=3696\", \"54\u00d791=4914\"],\n  [\"49\u00d782=4018\", \"98\u00d775=7350\"],\n  [\"84\u00d717=1428\", \"25\u00d717=425\"],\n  [\"59\u00d755=3245\", \"63\u00d799=6237\"],\n  [\"38\u00d746=1748\", \"64\u00d760=3840\"],\n  [\"86\u00d772=6192\", \"92\u00d754=4968\"],\n  [\"16\u00d748=768\", \"23\u00d787=2001\"],\n  [\"11\u00d789=979\", \"35\u00d744=1540\"],\n  [\"17\u00d768=1156\", \"35\u00d790=3150\"],\n  [\"41\u00d745=1845\", \"28\u00d761=1708\"],\n  [\"75\u00d764=4800\", \"85\u00d726=2210\"],\n  [\"48\u00d754=2592\", \"41\u00d797=3977\"],\n  [\"93\u00d741=3813\", \"58\u00d769=4002\"],\n  [\"61\u00d722=1342\", \"39\u00d735=1365\"],\n  [\"70\u00d765=4550\", \"77\u00d781=6237\"],\n  [\"13\u00d788=1144\", \"39\u00d729=1131\"],\n  [\"84\u00d789=7476\", \"96\u00d746=4416\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(`Text not found: ${oldText}`);\n  }\n\n  for (const item of results.items) {\n    item.insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Replace the date line and the 25 two-digit-by-two-digit multiplication\n# answers in the practice sheet with the new day's values. Every old\n# value below is unique in the document, so a simple exact-text\n# Find/Replace (one pair at a time) reproduces the diff exactly.\n\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"2025-03-27 Thursday\", \"2025-03-28 Friday\"),\n    @(\"56\u00d749=2744\", \"71\u00d754=3834\"),\n    @(\"53\u00d737=1961\", \"77\u00d746=3542\"),\n    @(\"13\u00d742=546\", \"76\u00d788=6688\"),\n    @(\"71\u00d770=4970\", \"96\u00d738=3648\"),\n    @(\"96\u00d719=1824\", \"30\u00d764=1920\"),\n    @(\"60\u00d777=4620\", \"58\u00d731=1798\"),\n    @(\"19\u00d782=1558\", \"66\u00d732=2112\"),\n    @(\"46\u00d794=4324\", \"12\u00d788=1056\"),\n    @(\"42\u00d788=3696\", \"54\u00d791=4914\"),\n    @(\"49\u00d782=4018\", \"98\u00d775=7350\"),\n    @(\"84\u00d717=1428\", \"25\u00d717=425\"),\n    @(\"59\u00d755=3245\", \"63\u00d799=6237\"),\n    @(\"38\u00d746=1748\", \"64\u00d760=3840\"),\n    @(\"86\u00d772=6192\", \"92\u00d754=4968\"),\n    @(\"16\u00d748=768\", \"23\u00d787=2001\"),\n    @(\"11\u00d789=979\", \"35\u00d744=1540\"),\n    @(\"17\u00d768=1156\", \"35\u00d790=3150\"),\n    @(\"41\u00d745=1845\", \"28\u00d761=1708\"),\n    @(\"75\u00d764=4800\", \"85\u00d726=2210\"),\n    @(\"48\u00d754=2592\", \"41\u00d797=3977\"),\n    @(\"93\u00d741=3813\", \"58\u00d769=4002\"),\n    @(\"61\u00d722=1342\", \"39\u00d735=1365\"),\n    @(\"70\u00d765=4550\", \"77\u00d781=6237\"),\n    @(\"13\u00d788=1144\", \"39\u00d729=1131\"),\n    @(\"84\u00d789=7476\", \"96\u00d746=4416\")\n)\n\nforeach ($pair in $replacements) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n\n    $rng = $d.Content\n    $rng.Find.ClearFormatting()\n    $found = $rng.Find.Execute($oldText, $false, $false, $false, $false, $false, $true, 1, $false, $newText, 2)\n    if (-not $found) {\n        throw \"Text not found: $oldText\"\n    }\n}\n"}
